$d = $word.ActiveDocument

# The "connect;solid" text run had a trailing ";solid" segment appended
# as its own run. Locate that exact run's text via Find (which narrows
# the Range to the match, same as real Word COM) and delete it outright,
# leaving the neighboring "connect" run untouched.
$rng = $d.Content
$found = $rng.Find.Execute(";solid")
if ($found) {
    $rng.Delete()
}
